$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the header style (bold, bordered, centered) from an existing
# header cell (AA1) onto the three new header cells so they match the
# rest of row 1 exactly.
$ws.Range("AA1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (1999 season: 74-88-0) repeated for every player row.
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
